# Apply the recorded changes to rows 63-66 of the "Artfynd" sheet.
# The underlying edit re-shuffles three observation records (previously
# stored in rows 63, 64 and 66) and rounds the Ost/Nord (Q/R) coordinates
# of all four affected rows (63-66) to whole metres.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 63  (becomes the former row-64 record: Svavelriska)
# ---------------------------------------------------------------
$ws.Range("A63").Value = 111683850
$ws.Range("B63").Value = 90332
$ws.Range("E63").Value = 4769
$ws.Range("F63").Value = "Svavelriska"
$ws.Range("G63").Value = "Lactarius scrobiculatus"
$ws.Range("H63").Value = "(Scop.:Fr.) Fr."

# "Antal" (I) is stored as text in this workbook, so force text format
# before writing the numeric-looking value "3".
$ws.Range("I63").NumberFormat = "@"
$ws.Range("I63").Value = "3"
$ws.Range("J63").Value = "fruktkroppar"

$ws.Range("P63").Value = "Bergaholm, Tyresö kn, Srm"
$ws.Range("Q63").Value = 689075
$ws.Range("R63").Value = 6570320
$ws.Range("S63").Value = 20
$ws.Range("Z63").Value = "09:25"
$ws.Range("AB63").Value = "09:25"

# ---------------------------------------------------------------
# Row 64  (becomes the former row-66 record: Fjällig taggsvamp s.str.)
# ---------------------------------------------------------------
$ws.Range("A64").Value = 111683853
$ws.Range("B64").Value = 90687
$ws.Range("E64").Value = 5964
$ws.Range("F64").Value = "Fjällig taggsvamp s.str."
$ws.Range("G64").Value = "Sarcodon imbricatus s.str."
$ws.Range("H64").Value = "(L.:Fr.) P.Karst."

# This record carries no "Antal"/"Enhet" -> clear those two cells.
$ws.Range("I64").Value = ""
$ws.Range("J64").Value = ""

$ws.Range("P64").Value = "Fiskarsundet, Srm"
$ws.Range("Q64").Value = 689112
$ws.Range("R64").Value = 6570306
$ws.Range("S64").Value = 23
$ws.Range("Z64").Value = "09:34"
$ws.Range("AB64").Value = "09:34"

# ---------------------------------------------------------------
# Row 65 (Sanicula europaea record stays, only coordinates rounded)
# ---------------------------------------------------------------
$ws.Range("Q65").Value = 689112
$ws.Range("R65").Value = 6570306

# ---------------------------------------------------------------
# Row 66  (becomes the former row-63 record: Rödgul trumpetsvamp)
# ---------------------------------------------------------------
$ws.Range("A66").Value = 111683845
$ws.Range("B66").Value = 89183
$ws.Range("E66").Value = 3215
$ws.Range("F66").Value = "Rödgul trumpetsvamp"
$ws.Range("G66").Value = "Craterellus lutescens"
$ws.Range("H66").Value = "(Fr.) Fr."

$ws.Range("Q66").Value = 689112
$ws.Range("R66").Value = 6570306
$ws.Range("Z66").Value = "09:36"
$ws.Range("AB66").Value = "09:36"
